# Commit: "data for new anr-h2 systems"
# Adds five new ANR-H2 (HTSE-coupled) plant-type rows (23-27) to the
# NewTechFramework sheet, mirroring the existing ANRElec rows (18-22)
# but tagged with the new ANRH2 DataSource/PlantCategory.
#
# Values are written one full column at a time (top to bottom) rather
# than row by row, so that brand-new shared-string entries land in the
# same left-to-right, column-major order Excel used originally.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewTechFramework")

$firstRow = 23
$plantTypes = @("iPWRHTSE", "HTGRHTSE", "PBRHTGRHTSE", "iMSRHTSE", "MicroHTSE")
$capacityMW = @(77, 164, 80, 141, 6.7)
$lifetimeYears = @(60, 60, 60, 60, 20)
$n = $plantTypes.Count

# Column -> either a single repeated value, or a per-row array.
$columns = [ordered]@{
    "A" = $plantTypes          # PlantType
    "B" = "ANRH2"               # DataSource
    "C" = "NA"                  # ATBTechnologyType
    "D" = "Nuclear Fuel"        # FuelType
    "E" = "h2"                  # ThermalOrRenewableOrStorage
    "F" = $capacityMW           # Capacity (MW)
    "K" = "Yes"                 # NSPSCompliant
    "L" = 0                     # NOxEmRate(lb/MMBtu)
    "M" = 0                     # SO2EmRate(lb/MMBtu)
    "N" = 0                     # CO2EmRate(lb/MMBtu)
    "O" = $lifetimeYears        # Lifetime(years)
    "Q" = "NA"                  # SO2 Scrubber
    "R" = "NA"                  # CoalType
    "S" = "NA"                  # Efficiency
    "U" = "NA"                  # Minimum Energy Capacity (MWh)
    "V" = "NA"                  # Maximum Charge Rate (MW)
    "W" = "NA"                  # ECAPEX(2012$/MWH)
    "X" = "ANRH2"                # PlantCategory
}

foreach ($col in $columns.Keys) {
    $value = $columns[$col]
    for ($i = 0; $i -lt $n; $i++) {
        $r = $firstRow + $i
        if ($value -is [array]) {
            $ws.Range("$col$r").Value = $value[$i]
        } else {
            $ws.Range("$col$r").Value = $value
        }
    }
}

$ws.Range("W31").Select()
